$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text replacements - these stay shared-string text on their own
# since Excel doesn't try to re-interpret them as numbers/dates.
$ws.Range("A2").Value = "shopping with Rubens"
$ws.Range("E2").Value = "Practical Steel Clock"
$ws.Range("F2").Value = "Blackberries"

# B2 ("1988-05-24") and G2 ("5.0") look like a date / a number, so a plain
# .Value assignment would make Excel silently convert them away from text.
# Mark the cells as Text first so the literal string is preserved.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1988-05-24"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5.0"

# H2 already reads "5.0" and keeps reading "5.0" after the edit (only the
# underlying shared-string slot it points at changes upstream) - so it is
# intentionally left untouched, preserving its original style/type.
